$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting so numeric-looking price
# strings (e.g. "1.006") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.879.93"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "1.831.35"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "310.55"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.4611"
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("D8").Value = "0.3671"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("D9").Value = "0.07168"
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("D10").Value = "0.8767"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").Value = "0.07893"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "19.58"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").Value = "1.855.06"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "5.336"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "6.380"
$ws.Range("E15").Value = "  -3.37%  "
$ws.Range("D16").Value = "87.50"
$ws.Range("E16").Value = "  -5.33%  "
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "0.000008725"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "26.919.94"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "14.45"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("D22").Value = "4.996"
$ws.Range("E22").Value = "  -2.98%  "
$ws.Range("D23").Value = "10.44"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").Value = "1.985"
$ws.Range("E24").Value = "  +4.56%  "
$ws.Range("D25").Value = "150.73"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("D26").Value = "18.22"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("D27").Value = "1.966"
$ws.Range("E27").Value = "  -5.38%  "
$ws.Range("D28").Value = "113.46"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("E29").Value = "  -4.10%  "
$ws.Range("D30").Value = "0.08839"
$ws.Range("D31").Value = "3.128"
$ws.Range("E31").Value = "  +3.26%  "
$ws.Range("D32").Value = "0.7541"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "4.454"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").Value = "1.128"
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("D35").Value = "2.582"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").Value = "1.087"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").Value = "0.01936"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("D38").Value = "2.927"
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("D39").Value = "0.05127"
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("D40").Value = "6.902"
$ws.Range("E40").Value = "  -3.46%  "
$ws.Range("D41").Value = "0.4977"
$ws.Range("E41").Value = "  -4.12%  "
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").Value = "8.326"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "0.4674"
$ws.Range("E44").Value = "  -3.84%  "
$ws.Range("D45").Value = "1.006"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "10.13"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").Value = "102.23"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("D50").Value = "64.50"
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("D51").Value = "36.34"
$ws.Range("E51").Value = "  -2.36%  "
